$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "sahing" -> "saying" in the "So when we say that the kingdom..."
# paragraph. The canonical edit splits the run into three pieces ("...sa",
# "y", "ing...") rather than collapsing to one run, so we reproduce that by
# inserting the "y" and deleting the stray "h" via small Range operations,
# then nudging a (no-op) Bold toggle across the new run boundaries so the
# engine keeps them as separate runs with clean/empty run-properties.
# ---------------------------------------------------------------------------

$findRange = $d.Content
$found = $findRange.Find.Execute("sahing", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $findRange.Start

    # locate the enclosing paragraph so we can re-touch the untouched head
    # and tail runs after the insert/delete below
    $paraRange = $d.Range($start, $start)
    $para = $paraRange.Paragraphs(1)
    $pStart = $para.Range.Start
    $pEnd = $para.Range.End - 1   # exclude the paragraph mark

    # "sahing" = s(0) a(1) h(2) i(3) n(4) g(5) relative to $start
    # Insert "y" right before the "h"
    $hRange = $d.Range($start + 2, $start + 3)
    $hRange.InsertBefore("y")

    # Delete the stray "h" (now shifted one to the right, at $start+3)
    $hRange2 = $d.Range($start + 3, $start + 4)
    $hRange2.Delete()

    # Split out the inserted "y" into its own run (toggling a boolean back
    # to its original value forces a run split without changing formatting)
    $yRange = $d.Range($start + 2, $start + 3)
    $yRange.Bold = 1
    $yRange.Bold = 0

    # Re-stamp clean run-properties on the untouched head ("...same as sa")
    # and tail ("ing...judgment?") runs that lost their <w:rPr/> during the
    # InsertBefore/Delete above.
    $run1 = $d.Range($pStart, $start + 2)
    $run1.Bold = 1
    $run1.Bold = 0

    $run3 = $d.Range($start + 3, $pEnd)
    $run3.Bold = 1
    $run3.Bold = 0
}

# ---------------------------------------------------------------------------
# Change 2: remove the duplicated "Me:" / "AI:" Q&A block about the
# statement of faith teaching on the "kingdom will be consummated" (it was
# pasted in twice, once under a "Me:" heading and once mistakenly again
# under an "AI:" heading).
# ---------------------------------------------------------------------------

$target = 'Our statement of faith teaches that the "kingdom will be consummated" at the bodily return of christ as final judge.'

$count = $d.Paragraphs.Count
$firstMeIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($ptext -eq "Me:") {
        $next = $d.Paragraphs($i + 1)
        $ntext = $next.Range.Text.TrimEnd([char]13, [char]7)
        if ($ntext.StartsWith($target)) {
            $firstMeIndex = $i
            break
        }
    }
}

if ($firstMeIndex -gt 0) {
    # The duplicated block is: Me: / statement / (blank) / AI: / statement / (blank)
    # i.e. 6 paragraphs starting at $firstMeIndex, followed by the next
    # surviving "Me:" heading paragraph.
    $blockStartPara = $d.Paragraphs($firstMeIndex)
    $afterBlockPara = $d.Paragraphs($firstMeIndex + 6)

    $delStart = $blockStartPara.Range.Start
    $delEnd = $afterBlockPara.Range.Start

    $delRange = $d.Range($delStart, $delEnd)
    $delRange.Delete()
}
